$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9993789792060852
$ws.Range("B1").Value = 4.112075805664062
$ws.Range("C1").Value = 2.346004009246826
$ws.Range("D1").Value = 1.71368396282196
$ws.Range("E1").Value = 1.345895648002625
